$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New daily data rows (302:328), continuing the series through 2021-07-25
# (Excel serial dates 44376-44402), columns: A=date, B=nuovi pos.,
# C=somma mobile 7gg., D=somma mobile 7gg. per 100mila abitanti.
$data = @(
    @(44376, 2, 7, 38.93647791745467),
    @(44377, 0, 5, 27.81176994103905),
    @(44378, 0, 5, 27.81176994103905),
    @(44379, 0, 4, 22.24941595283124),
    @(44380, 1, 3, 16.68706196462343),
    @(44381, 0, 3, 16.68706196462343),
    @(44382, 0, 3, 16.68706196462343),
    @(44383, 0, 1, 5.56235398820781),
    @(44384, 0, 1, 5.56235398820781),
    @(44385, 0, 1, 5.56235398820781),
    @(44386, 1, 2, 11.12470797641562),
    @(44387, 0, 1, 5.56235398820781),
    @(44388, 0, 1, 5.56235398820781),
    @(44389, 0, 1, 5.56235398820781),
    @(44390, 0, 1, 5.56235398820781),
    @(44391, 0, 1, 5.56235398820781),
    @(44392, 0, 1, 5.56235398820781),
    @(44393, 0, 0, 0),
    @(44394, 0, 0, 0),
    @(44395, 1, 1, 5.56235398820781),
    @(44396, 0, 1, 5.56235398820781),
    @(44397, 0, 1, 5.56235398820781),
    @(44398, 0, 1, 5.56235398820781),
    @(44399, 0, 1, 5.56235398820781),
    @(44400, 0, 1, 5.56235398820781),
    @(44401, 1, 2, 11.12470797641562),
    @(44402, 2, 3, 16.68706196462343)
)

$firstRow = 302
$lastRow = $firstRow + $data.Count - 1

# Extend column A formatting (the "YYYY-MM-DD HH:MM:SS" date style used by
# the rest of the column, style index s="2") down into the new rows before
# writing values, so the new cells keep matching formatting.
$ws.Range("A301").Copy($ws.Range("A$firstRow`:A$lastRow"))

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $firstRow + $i
    $vals = $data[$i]
    $ws.Cells.Item($r, 1).Value = $vals[0]
    $ws.Cells.Item($r, 2).Value = $vals[1]
    $ws.Cells.Item($r, 3).Value = $vals[2]
    $ws.Cells.Item($r, 4).Value = $vals[3]
}

Write-Output "Added rows $firstRow to $lastRow"